$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out rows beyond the new data extent (old sheet had up to row 110; new data ends at row 101)
$ws.Range("A102:B110").ClearContents()

$names = @(
  "3535 Opal Meadow Heights Aged Care Community Meadow Heights",
  "95 Napier Street Apartment Complex Fitzroy",
  "Al Haj Halal Meats Glenroy",
  "Al-Taqwa College Truganina",
  "Amiga Montessori Craigieburn",
  "Australia Post Distribution Centre Sunshine West",
  "Baxter Foods Australia Campbellfield",
  "Budget Car and Truck Rentals Campbellfield",
  "CS Square Caroline Springs",
  "Cafe Roco Dandenong",
  "Campbellfield Ford Complex Vaccination Clinic Campbellfield",
  "Cannie Road Construction Site Cannie",
  "Caroline Springs Police Station",
  "Cedars Medical Clinic Coburg",
  "Chemist Warehouse Campbellfield DC",
  "Chemist Warehouse Fillo Drive Somerton",
  "City of Wyndham Community",
  "Coles Broadmeadows Central Shopping Centre",
  "Coles Campbellfield Plaza Campbellfield",
  "Coles Coburg North Village",
  "Coles Pakenham Place Shopping Centre",
  "Coles Roxburgh Village Roxburgh Park",
  "Community Kids Bayswater Early Education Centre Bayswater North",
  "Community Kids Meadow Heights",
  "Construction Site 1 Warde Street Footscray",
  "Construction Site Olea Apartment Caulfield North",
  "Costco Wholesale Epping",
  "Crusader Caravans Epping",
  "DayHab Rehabilitation Treatment Centre Ringwood East",
  "Direct Freight Express Cambellfield",
  "Disability Residence Life without Barriers Ashwood",
  "Don Watson Coldstore Derrimut",
  "Epworth Healthcare Epworth Richmond Emergency Department",
  "FedEx Station Melbourne Airport",
  "Fine Food Holdings Pty Ltd Dandenong South",
  "Fitzroy Community School Fitzroy North",
  "Fonterra Manufacturing Workplace Campbellfield",
  "General Foods Campbellfield",
  "Gladstone Parade Early Learning & Kinder Glenroy",
  "Goodstart Early Learning Altona",
  "Green Leaves Early Learning Cairnlea",
  "Green Leaves Early Learning Centre Highlands Craigieburn",
  "Hamilton Marino 236 Jasper Road McKinnon",
  "Hello Fresh Warehouse Ravenhall",
  "Hickory Construction Site Chadstone Car Park Malvern East",
  "IGA Meadow Heights Shopping Centre Meadow Heights",
  "ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine",
  "Ibis Kingsgate Hotel Melbourne",
  "Ilim Learning Sanctuary Glenroy",
  "Industrial Galvanizers Valmont Coatings Campbellfield",
  "Inghams Enterprises Thomastown",
  "KFC Fawkner",
  "Kasr Sweets Coolaroo",
  "Kids House Early Learning Cheltenham",
  "Kippers Seafood Werribee",
  "Kool Kidz Childcare Narre Warren",
  "Level Crossing Removal Project Lilydale Construction Site John Street",
  "Lineage Logistics Laverton North",
  "Linfox Somerton National Distribution Centre Somerton",
  "McDonald's Craigieburn North",
  "Mecca D.C Warehouse Melbourne Airport",
  "Melbourne Assessment Prison West Melbourne",
  "Melbourne Metropolitan Remand Centre Ravenhall",
  "Melbourne West Police Station Docklands",
  "Mill Park Police Station Mill Park",
  "MyCentre Childcare Broadmeadows",
  "National Gallery of Victoria Melbourne",
  "Nido Early School Ascot Vale",
  "Nido Early School Glenroy",
  "Nido Early School Moonee Ponds",
  "Northern Health Northern Hospital Epping Emergency Department Tier 1B",
  "Northern Health The Northern Hospital Epping",
  "OnQ Plumbing and Excavations Craigieburn",
  "Oporto Coolaroo",
  "Oscar Romero Catholic Primary School Craigieburn",
  "Our Lady Help of Christian's Primary School Brunswick East",
  "Pacific Meat Thomastown",
  "Panorama Construction Site Whitehorse Rd Box Hill",
  "Private Residence Northern Community Services Fawkner",
  "Ramsay Health Care Warringal Private Hospital Heidelberg",
  "Ravenhall Correctional Centre Ravenhall",
  "Richmond Quarter 261-271 Bridge Road Construction Site Richmond",
  "Sacca's Fruit World Broadmeadows Central Shopping Centre",
  "Salta Drive Construction Site Rangedale Drainage Altona North",
  "St Margaret's Primary School OSHC Maribyrnong",
  "St Vincents Hospital Emergency Department Melbourne",
  "Sultan Halal Meats & Poultry Campbellfield",
  "Tek Foods Somerton",
  "The Huntly-Goornong Rail Works",
  "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B",
  "The Royal Melbourne Hospital AMU Ward Parkville",
  "ThorwestenCabinets Pakenham",
  "Total Window Concepts Hoppers Crossing",
  "Unilodge College Square Student Accommodation 570 Lygon Street Carlton",
  "Wallaby Childcare Wollert",
  "Werribee Mercy Hospital Emergency Department",
  "Western Health Footscray Hospital Emergency Department",
  "Western Health Sunshine Hospital Emergency Department",
  "Woodlands Long Day Care and Kindergarten Roxburgh Park",
  "Yara Childcare Centre Truganina"
)

$values = @(
  27,
  5,
  23,
  9,
  8,
  5,
  5,
  5,
  14,
  5,
  9,
  5,
  10,
  17,
  6,
  11,
  6,
  5,
  9,
  25,
  10,
  5,
  12,
  7,
  5,
  16,
  29,
  23,
  5,
  11,
  5,
  5,
  6,
  13,
  9,
  9,
  10,
  13,
  6,
  11,
  5,
  16,
  11,
  5,
  5,
  6,
  11,
  5,
  6,
  19,
  6,
  5,
  5,
  8,
  6,
  9,
  7,
  9,
  10,
  5,
  7,
  7,
  8,
  6,
  7,
  13,
  9,
  7,
  12,
  11,
  62,
  13,
  14,
  6,
  5,
  10,
  5,
  8,
  5,
  7,
  7,
  11,
  6,
  7,
  11,
  10,
  5,
  13,
  6,
  17,
  25,
  14,
  6,
  8,
  13,
  11,
  8,
  9,
  5,
  7
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

Write-Host "Done. Rows written: $($names.Count)"